$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.893.66'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.358.53'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.28'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.668'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.50%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.58'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.90%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.608'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '60.74'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +6.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '33.81'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.108'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.26'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.22'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.910'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.354.27'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.914.07'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '77.71'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.54'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '253.23'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.86'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.72%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.45'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.29%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '176.07'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.28'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.95%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.53%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0746'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.06'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.34'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.79'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.53'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.41'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0276'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.48'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +16.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '65.37'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +12.91%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.67%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.107'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -6.54%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.07'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.200'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.03%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.20%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.44'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.16'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '98.17'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.87'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.83%  '
